# Update "想去人数" (column F) values on the "展览" (sheet1) and
# "全部类型" (sheet4) worksheets to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 466
$ws1.Range("F6").Value = 279
$ws1.Range("F7").Value = 2542
$ws1.Range("F8").Value = 433
$ws1.Range("F9").Value = 6845
$ws1.Range("F10").Value = 185
$ws1.Range("F11").Value = 434
$ws1.Range("F12").Value = 3
$ws1.Range("F13").Value = 33

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 466
$ws4.Range("F6").Value = 279
$ws4.Range("F9").Value = 2542
$ws4.Range("F10").Value = 433
$ws4.Range("F11").Value = 6845
$ws4.Range("F12").Value = 185
$ws4.Range("F13").Value = 434
$ws4.Range("F14").Value = 3
$ws4.Range("F17").Value = 33
